# Update the "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to reflect the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 301
$wsExpo.Range("F3").Value = 1203
$wsExpo.Range("F4").Value = 16844
$wsExpo.Range("F5").Value = 31
$wsExpo.Range("F7").Value = 68
$wsExpo.Range("F9").Value = 384
$wsExpo.Range("F10").Value = 222
$wsExpo.Range("F12").Value = 11674
$wsExpo.Range("F14").Value = 1350
$wsExpo.Range("F15").Value = 4634
$wsExpo.Range("F16").Value = 455
$wsExpo.Range("F17").Value = 405
$wsExpo.Range("F20").Value = 339

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 301
$wsAll.Range("F4").Value = 1203
$wsAll.Range("F5").Value = 16844
$wsAll.Range("F6").Value = 31
$wsAll.Range("F8").Value = 68
$wsAll.Range("F10").Value = 384
$wsAll.Range("F11").Value = 222
$wsAll.Range("F15").Value = 11674
$wsAll.Range("F17").Value = 1350
$wsAll.Range("F18").Value = 4634
$wsAll.Range("F19").Value = 455
$wsAll.Range("F20").Value = 405
$wsAll.Range("F23").Value = 339
